$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Caso1")
$values = @(
    1.003542304039001,
    1.002823948860168,
    1.000772953033447,
    1.009588599205017,
    1.011250615119934,
    1.000048756599426,
    1.002289295196533,
    1.000875473022461,
    1.002793312072754,
    1.004420518875122,
    0.9982147812843323,
    0.9986172914505005,
    1.004990696907043,
    1.008272647857666,
    1.007127642631531,
    0.9997249841690063,
    1.000653862953186,
    1.004794001579285
)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $values[$i]
}

$ws = $wb.Worksheets.Item("Caso2")
$values = @(
    1.002911329269409,
    1.002636551856995,
    1.000641822814941,
    1.009082317352295,
    1.010901570320129,
    0.9995949268341064,
    1.001675724983215,
    1.000414371490479,
    1.002419471740723,
    1.004103422164917,
    0.9981140494346619,
    0.997689962387085,
    1.004761219024658,
    1.008124828338623,
    1.006569266319275,
    0.9994498491287231,
    1.000728845596313,
    1.004410862922668
)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $values[$i]
}

$ws = $wb.Worksheets.Item("Caso3")
$values = @(
    1.002438187599182,
    1.001886963844299,
    0.9998671412467957,
    1.008527159690857,
    1.010313749313354,
    0.9992058277130127,
    1.001111388206482,
    0.9999788403511047,
    1.001765131950378,
    1.003379106521606,
    0.9974510073661804,
    0.9976387619972229,
    1.00419819355011,
    1.007380127906799,
    1.006126046180725,
    0.9987671971321106,
    0.9996189475059509,
    1.003853678703308
)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $values[$i]
}

